# Commit: feat: add 2022-Q4 data
#
# Semantically this edit:
#   1. Inserts a new worksheet named "2022-Q4" right after "总计" (so it
#      becomes the 2nd tab), containing the quarterly fund-holdings detail
#      for the same 4 funds that appear in "2022-Q3", but with the new
#      quarter's numbers.
#   2. Adds a corresponding summary row for "2022-Q4" at the top of the
#      data in the "总计" sheet (pushing the existing history down by one
#      row).
#
# All of the other worksheets keep their original data untouched - the
# fact that their underlying xl/worksheets/sheetN.xml part numbers shift
# in the saved package is just a natural side effect of inserting a new
# sheet before them; the worksheets' names/content are unchanged.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new row 2 for 2022-Q4.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows(2).Insert()

# Copy the formatting from the row below (the old row 2, now row 3) onto
# the freshly inserted row so the styling (bold index column, etc.)
# matches the rest of the table exactly.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 4
$summary.Cells.Item(2, 4).Value = 1.28

# ---------------------------------------------------------------------
# 2. New "2022-Q4" worksheet: duplicate "2022-Q3" (same layout/styling,
#    same 4 funds in the same order) immediately before it, rename, and
#    overwrite the quarter-specific numbers.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# D:G hold numbers formatted as text in this workbook (t="inlineStr" in
# the original OOXML) - force text storage so "6.89" etc. aren't
# reinterpreted as numeric values.
$q4.Range("D2:G5").NumberFormat = "@"

# Row 2 - 002207 前海开源金银珠宝主题精选混合C
$q4.Cells.Item(2, 4).Value = "6.89"
$q4.Cells.Item(2, 5).Value = "90.36"
$q4.Cells.Item(2, 6).Value = "7.57"
$q4.Cells.Item(2, 7).Value = "0.5216"
$q4.Cells.Item(2, 8).Value = 9

# Row 3 - 001302 前海开源金银珠宝主题精选混合A
$q4.Cells.Item(3, 4).Value = "3.55"
$q4.Cells.Item(3, 5).Value = "90.36"
$q4.Cells.Item(3, 6).Value = "7.57"
$q4.Cells.Item(3, 7).Value = "0.2687"
$q4.Cells.Item(3, 8).Value = 9

# Row 4 - 003304 前海开源沪港深核心资源灵活配置混合A
$q4.Cells.Item(4, 4).Value = "3.30"
$q4.Cells.Item(4, 5).Value = "90.48"
$q4.Cells.Item(4, 6).Value = "7.56"
$q4.Cells.Item(4, 7).Value = "0.2495"
$q4.Cells.Item(4, 8).Value = 7

# Row 5 - 003305 前海开源沪港深核心资源灵活配置混合C
$q4.Cells.Item(5, 4).Value = "3.17"
$q4.Cells.Item(5, 5).Value = "90.48"
$q4.Cells.Item(5, 6).Value = "7.56"
$q4.Cells.Item(5, 7).Value = "0.2397"
$q4.Cells.Item(5, 8).Value = 7

# ---------------------------------------------------------------------
# 3. Restore the originally active sheet ("总计" / tab 0) so the new
#    sheet doesn't steal the workbook's active-tab selection.
# ---------------------------------------------------------------------
$summary.Activate()
